$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Update Runmode for TestCase_F2 (row 3) from N to Y
$ws.Range("C3").Value = "Y"

# Mark Results as SKIP for TestCase_F2, F3, F4 (rows 3-5); F1 (row 2) stays PASS
$ws.Range("D3").Value = "SKIP"
$ws.Range("D4").Value = "SKIP"
$ws.Range("D5").Value = "SKIP"

# Update the active selection to C4
$ws.Activate()
$ws.Range("C4").Select()
